$d = $word.ActiveDocument

# Locate the existing "32. The Gerontologist" entry so we can append the
# new "33. Annals of Internal Medicine" bibliography entry right after it.
$hit = $d.Content
$found = $hit.Find.Execute("The Gerontologist", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate 'The Gerontologist' paragraph"
}

$sourcePara = $hit.Paragraphs(1)

# Insert a brand-new paragraph right after it; Word clones the paragraph
# formatting (firstLine indent) of the paragraph it was split from.
$sourcePara.Range.InsertParagraphAfter() | Out-Null
$newPara = $sourcePara.Next()
$target = $newPara.Range

# Build the two runs ("33. " and "Annals of Internal Medicine") explicitly
# via OOXML so they remain separate <w:r> elements, matching how the
# entry was originally authored, and keep the firstLine indent.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' +
       '<w:r><w:t xml:space="preserve">33. </w:t></w:r>' +
       '<w:r><w:t>Annals of Internal Medicine</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml) | Out-Null
